$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, pushing existing rows 132-238 down to 133-239
$ws.Rows.Item(132).Insert()

# Populate the new row 132 with the latest weekly record
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 44634
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112008
$ws.Cells.Item(132, 7).Value = "Coliflor"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 200
$ws.Cells.Item(132, 11).Value = 1200
$ws.Cells.Item(132, 12).Value = 1300
$ws.Cells.Item(132, 13).Value = 1250
$ws.Cells.Item(132, 14).Value = "`$/unidad"
$ws.Cells.Item(132, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(132, 16).Value = 1250
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = "Hortaliza"
